$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8: date 8/2/2021 in column A (matching the d-mmm style already used
# by the other date rows), and a new shared-string question in column B.
$ws.Cells.Item(8, 1).Value = [datetime]"2021-08-02"
$ws.Cells.Item(8, 1).NumberFormat = $ws.Cells.Item(6, 1).NumberFormat
$ws.Cells.Item(8, 2).Value = "415_AddStrings"

# New row 9: a second question added the same day, no date in column A.
$ws.Cells.Item(9, 2).Value = "67_AddBinary"

# Move the selection to where the author left off editing.
$ws.Range("L14").Select()
